# Task 555 - added validation for gross_annual_revenue and gross_annual_revenue_flag fields
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# New header columns: X = gross_annual_revenue_flag, Y = gross_annual_revenue
$ws.Range("X1").Value = "gross_annual_revenue_flag"
$ws.Range("Y1").Value = "gross_annual_revenue"

# Row 2
$ws.Range("X2").Value = 900
$ws.Range("Y2").Value = 100000

# Row 3
$ws.Range("X3").Value = 988

# Row 4
$ws.Range("X4").Value = 900

# Row 5
$ws.Range("X5").Value = 988
$ws.Range("Y5").Value = 300000

# Row 6
$ws.Range("Y6").Value = 200000

# Row 7
$ws.Range("X7").Value = 999

# Row 8
$ws.Range("X8").Value = 990
$ws.Range("Y8").Value = 50000

# Row 9
$ws.Range("X9").Value = 900
$ws.Range("Y9").Value = 45000

# Row 10
$ws.Range("X10").Value = 988

# Row 11
$ws.Range("X11").Value = 900
$ws.Range("Y11").Value = 50000

# Update the active selection to match the edited workbook state
$ws.Range("Y13").Select() | Out-Null
